# Add a new row 50 (image-processing leak data point) to each of the four
# FE_* worksheets, mirroring the existing row 48/49 layout.

$wb = $excel.ActiveWorkbook

$newRow = 50

# Per-sheet values for the new row. Column order: A time, B total-len hex,
# C id hex, D actual-len hex, E checksum hex, F total-len dec, G id dec,
# H actual-len dec, I checksum dec.
$rowsData = @{
    "FE_LFT_#1" = @{
        A = 45836.49679398148
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x58"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 344
        I = 15
    }
    "FE_LFT_#2" = @{
        A = 45836.49679398148
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x6C"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 364
        I = 14
    }
    "FE_PLT_#1" = @{
        A = 45836.49679398148
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x69"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 105
        I = 3
    }
    "FE_PLT_#2" = @{
        A = 45836.49679398148
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x68"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 104
        I = 3
    }
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowsData[$ws.Name]
    if ($data -eq $null) {
        continue
    }

    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}

Write-Output "done"
